$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 96, shifting existing rows 96:185 down to 97:186
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with its data
$ws.Cells.Item(96, 1).Value = 4
$ws.Cells.Item(96, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(96, 3).Value = "Los Lagos"
$ws.Cells.Item(96, 4).Value = 44586
$ws.Cells.Item(96, 5).Value = 10
$ws.Cells.Item(96, 6).Value = 100112032
$ws.Cells.Item(96, 7).Value = "Zapallo italiano"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 200
$ws.Cells.Item(96, 11).Value = 16000
$ws.Cells.Item(96, 12).Value = 17000
$ws.Cells.Item(96, 13).Value = 16500
$ws.Cells.Item(96, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(96, 15).Value = "Región Metropolitana"
$ws.Cells.Item(96, 16).Value = 330
$ws.Cells.Item(96, 17).Value = 50
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# Apply the same date style as the other rows in column D
$ws.Cells.Item(96, 4).NumberFormat = $ws.Cells.Item(97, 4).NumberFormat
